# Edit workbook: add "2022-Q4" sheet with fund holding data, and update
# the "总计" (summary) sheet with a new leading row for 2022-Q4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($q3)
$ws.Name = "2022-Q4"

# Helper: force a cell's value to be written as *text*, even when the
# string looks like a number (e.g. fund codes such as "012284" or
# percentage-like values such as "30.82"), matching the source data
# which stores these as text rather than numbers.
function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 2. Header row (row 1), columns B:H - copy the bold/bordered/centered
#    style already used for header cells elsewhere in the workbook.
# ---------------------------------------------------------------------
$summary.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Copy the bold/bordered/centered style used for the "A" index
#    column onto A2:A24 of the new sheet.
# ---------------------------------------------------------------------
$summary.Range("A2").Copy()
$ws.Range("A2:A24").PasteSpecial(-4122)

$data = @(
    @(0, "163807", "中银优选混合A", "30.82", "70.19", "2.13", "0.6565", 5, $true),
    @(1, "012284", "光大保德信健康优加混合", "13.09", "87.41", "3.44", "0.4503", 8, $true),
    @(2, "012631", "中银优选混合C", "12.31", "70.19", "2.13", "0.2622", 5, $true),
    @(3, "360005", "光大保德信红利混合", "3.84", "79.76", "3.44", "0.1321", 5, $true),
    @(4, "010090", "中信建投医药健康混合A", "2.95", "95.01", "3.58", "0.1056", 10, $true),
    @(5, "009414", "中银大健康股票A", "2.85", "88.80", "2.63", "0.0750", 8, $true),
    @(6, "010091", "中信建投医药健康混合C", "1.98", "95.01", "3.58", "0.0709", 10, $true),
    @(7, "008905", "嘉合锦鹏添利混合A", "3.64", "22.86", "1.05", "0.0382", 3, $true),
    @(8, "008422", "中融研发创新混合A", "1.07", "61.65", "3.47", "0.0371", 7, $true),
    @(9, "004671", "中融核心成长灵活配置混合", "1.14", "65.86", "3.04", "0.0347", 9, $true),
    @(10, "002504", "鹏华金鼎灵活配置混合A", "0.60", "88.96", "5.20", "0.0312", 4, $true),
    @(11, "005293", "诺德新旺灵活配置混合", "0.55", "92.81", "5.28", "0.0290", 9, $true),
    @(12, "010009", "中融成长优选混合C", "1.05", "60.81", "2.65", "0.0278", 9, $true),
    @(13, "008906", "嘉合锦鹏添利混合C", "2.16", "22.86", "1.05", "0.0227", 3, $true),
    @(14, "008423", "中融研发创新混合C", "0.62", "61.65", "3.47", "0.0215", 7, $true),
    @(15, "163818", "中银中小盘成长混合", "0.77", "87.05", "2.71", "0.0209", 5, $true),
    @(16, "206013", "鹏华宏观灵活配置混合", "0.37", "72.26", "4.22", "0.0156", 6, $true),
    @(17, "010008", "中融成长优选混合A", "0.57", "60.81", "2.65", "0.0151", 9, $true),
    @(18, "000591", "中银健康生活混合", "0.49", "71.03", "2.51", "0.0123", 4, $true),
    @(19, "002505", "鹏华金鼎灵活配置混合C", "0.18", "88.96", "5.20", "0.0094", 4, $true),
    @(20, "005545", "中银改革红利灵活配置混合A", "0.46", "65.05", "2.03", "0.0093", 7, $true),
    @(21, "010321", "中银大健康股票C", "0.34", "88.80", "2.63", "0.0089", 8, $true),
    @(22, "016480", "中银改革红利灵活配置混合C", "0.00", "65.05", "2.03", 0, 7, $false),
)


for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $ws.Range("A$row").Value = $rec[0]
    Set-TextCell $ws.Range("B$row") $rec[1]
    Set-TextCell $ws.Range("C$row") $rec[2]
    Set-TextCell $ws.Range("D$row") $rec[3]
    Set-TextCell $ws.Range("E$row") $rec[4]
    Set-TextCell $ws.Range("F$row") $rec[5]

    if ($rec[8]) {
        Set-TextCell $ws.Range("G$row") $rec[6]
    } else {
        $ws.Range("G$row").Value = $rec[6]
    }

    $ws.Range("H$row").Value = $rec[7]
}

# ---------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert a new row for
#    2022-Q4 at the top of the data (row 2), pushing the existing
#    quarters down by one row.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 23
$summary.Range("D2").Value = 2.09

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
